$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.896.68"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.649.93"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.18"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.10"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.57"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.382"
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "3.122.13"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.71"
$ws.Range("E14").Value = "  +10.43%  "
$ws.Range("D15").Value = "60.872.77"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "2.667.48"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.61"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.67"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.535"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.97"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.13"
$ws.Range("E27").Value = "  +4.64%  "
$ws.Range("E28").Value = "  +6.67%  "
$ws.Range("D29").Value = "0.0₃0814"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("E30").Value = "  +6.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.37"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("E34").Value = "  +9.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  +5.79%  "
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "332.06"
$ws.Range("E38").Value = "  +12.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.02"
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.40"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.878"
$ws.Range("E41").Value = "  +4.32%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.55"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.67"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.615"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "2.111.64"
$ws.Range("E51").Value = "  +4.25%  "
